$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column L entirely (the "risk severe 20-29" computed column, =Kx/K7).
# This shifts the former M/N columns left into L/M for every row.
$ws.Range("L:L").Delete()

# Fix up the renamed English header labels in row 3 (columns that kept their
# position are simply re-labelled; L3/M3 are the ones that shifted left).
$ws.Range("B3").Value = "unvax"
$ws.Range("E3").Value = "unvax per 100k"
$ws.Range("G3").Value = "partial vax per 100k"
$ws.Range("L3").Value = "severe vax per 100k"
$ws.Range("M3").Value = "severe parital vax per 100k"

$ws.Range("A3").Select()
